# Updates the Price (D) and Volume(1h) (E) columns of the crypto
# tracker sheet to the freshly scraped values, per the
# "Updated symbol list" GitHub Actions commit.
#
# Values are written as literal text (matching the existing
# inline-string cells for these columns) rather than being left for
# Excel to auto-coerce into numbers/percentages, so "300.54" stays the
# text "300.54" and "-0.22%" stays the text "-0.22%" instead of turning
# into numeric/percentage cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row=2; D='300.54'; E='-0.22%' }
    @{ Row=3; D='31.82'; E='1.50%' }
    @{ Row=4; D=''; E='0.60%' }
    @{ Row=5; D='0.08157'; E='10.64%' }
    @{ Row=6; D='2.524'; E='3.88%' }
    @{ Row=7; D='7.847'; E='-1.38%' }
    @{ Row=8; D='3.869'; E='2.21%' }
    @{ Row=9; D='0.9252'; E='0.96%' }
    @{ Row=10; D='0.1758'; E='2.98%' }
    @{ Row=11; D='0.07414'; E='-2.60%' }
    @{ Row=12; D='0.08933'; E='10.19%' }
    @{ Row=13; D='0.03025'; E='-0.44%' }
    @{ Row=14; D='0.1002'; E='1.00%' }
    @{ Row=15; D='0.001521'; E='1.04%' }
    @{ Row=16; D='0.005997'; E='-2.75%' }
    @{ Row=17; D='3.606'; E='4.09%' }
    @{ Row=18; D=''; E='2.63%' }
    @{ Row=19; D=''; E='-1.02%' }
    @{ Row=20; D='0.1339'; E='0.28%' }
    @{ Row=21; D='4.077'; E='-12.24%' }
    @{ Row=22; D=''; E='7.30%' }
    @{ Row=23; D='0.04637'; E='-0.29%' }
    @{ Row=24; D='0.001247'; E='1.76%' }
    @{ Row=25; D='0.004544'; E='1.29%' }
    @{ Row=26; D='0.0001199'; E='-7.60%' }
    @{ Row=27; D='0.0003409'; E='82.15%' }
    @{ Row=39; D='0.01768'; E='1.84%' }
    @{ Row=40; D='0.04610'; E='1.93%' }
    @{ Row=41; D='0.006852'; E='-4.94%' }
    @{ Row=42; D='0.1380'; E='2.46%' }
    @{ Row=43; D='0.002208'; E='-0.80%' }
    @{ Row=44; D='0.009853'; E='-8.22%' }
    @{ Row=45; D='0.00006166'; E='-1.68%' }
    @{ Row=46; D='0.00000000750'; E='-0.03%' }
    @{ Row=47; D='0.8044'; E='-0.51%' }
    @{ Row=48; D='0.008398'; E='-15.95%' }
    @{ Row=49; D='0.00002100'; E='-0.03%' }
    @{ Row=50; D='0.0002000'; E='0.04%' }
)

foreach ($u in $updates) {
    $r = $u.Row

    if ($u.D -ne '') {
        $dCell = $ws.Range("D$r")
        $dCell.Value = "'" + $u.D
        $dCell.Style = "Normal"
    }

    $eCell = $ws.Range("E$r")
    $eCell.Value = "'" + $u.E
    $eCell.Style = "Normal"
}
